# Auto-generated Excel COM-interop script
# Applies the crypto-price/volume refresh to Sheet1's data rows,
# matching the upstream "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell($addr, $val) {
    # Price column (D) holds values such as "27.030.27" or "1.006" that
    # must stay literal text (not be re-interpreted as a number/date by
    # the usual Excel input parsing), exactly like the source data feed.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $val
}

Set-TextCell 'D2' '27.030.27'
$ws.Range('E2').Value = '  +0.10%  '
Set-TextCell 'D3' '1.829.83'
$ws.Range('E3').Value = '  +0.47%  '
Set-TextCell 'D4' '1.006'
$ws.Range('E4').Value = '  -0.53%  '
Set-TextCell 'D5' '311.74'
$ws.Range('E5').Value = '  +0.51%  '
Set-TextCell 'D6' '1.005'
$ws.Range('E6').Value = '  -0.45%  '
Set-TextCell 'D7' '0.4630'
$ws.Range('E7').Value = '  -0.27%  '
Set-TextCell 'D8' '0.3713'
$ws.Range('E8').Value = '  +2.16%  '
Set-TextCell 'D9' '0.07352'
$ws.Range('E9').Value = '  +0.81%  '
Set-TextCell 'D10' '0.8774'
$ws.Range('E10').Value = '  +1.36%  '
Set-TextCell 'D11' '0.07882'
$ws.Range('E11').Value = '  +4.15%  '
Set-TextCell 'D12' '19.79'
$ws.Range('E12').Value = '  -0.09%  '
Set-TextCell 'D13' '1.808.68'
$ws.Range('E13').Value = '  -1.18%  '
Set-TextCell 'D14' '5.340'
$ws.Range('E14').Value = '  -0.03%  '
Set-TextCell 'D15' '6.552'
$ws.Range('E15').Value = '  +1.36%  '
Set-TextCell 'D16' '91.51'
$ws.Range('E16').Value = '  -1.92%  '
Set-TextCell 'D17' '1.007'
$ws.Range('E17').Value = '  -0.34%  '
Set-TextCell 'D18' '0.000008862'
$ws.Range('E18').Value = '  +2.68%  '
Set-TextCell 'D19' '1.005'
$ws.Range('E19').Value = '  -0.51%  '
Set-TextCell 'D20' '14.80'
$ws.Range('E20').Value = '  +2.28%  '
Set-TextCell 'D21' '27.052.69'
$ws.Range('E21').Value = '  -0.75%  '
Set-TextCell 'D22' '5.111'
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('E23').Value = '  -0.25%  '
Set-TextCell 'D24' '2.023.35'
$ws.Range('E24').Value = '  -3.08%  '
Set-TextCell 'D25' '152.58'
$ws.Range('E25').Value = '  +0.57%  '
Set-TextCell 'D26' '1.849'
$ws.Range('E26').Value = '  -0.47%  '
Set-TextCell 'D27' '18.44'
$ws.Range('E27').Value = '  +1.04%  '
Set-TextCell 'D28' '2.041'
$ws.Range('E28').Value = '  -2.57%  '
Set-TextCell 'D29' '5.127'
$ws.Range('E29').Value = '  +1.10%  '
Set-TextCell 'D30' '115.71'
$ws.Range('E30').Value = '  -0.16%  '
Set-TextCell 'D31' '0.08897'
$ws.Range('E31').Value = '  -0.08%  '
Set-TextCell 'D32' '2.962'
$ws.Range('E32').Value = '  +0.27%  '
Set-TextCell 'D33' '0.7284'
$ws.Range('E33').Value = '  +0.01%  '
Set-TextCell 'D34' '4.449'
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('E35').Value = '  -0.25%  '
Set-TextCell 'D36' '2.468'
$ws.Range('E36').Value = '  -1.30%  '
Set-TextCell 'D37' '1.079'
$ws.Range('E37').Value = '  +0.48%  '
Set-TextCell 'D38' '0.01950'
$ws.Range('E38').Value = '  +1.78%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('E40').Value = '  +0.80%  '
Set-TextCell 'D41' '7.134'
$ws.Range('E41').Value = '  +0.31%  '
Set-TextCell 'D42' '0.5181'
$ws.Range('E42').Value = '  -0.42%  '
Set-TextCell 'D43' '0.1627'
$ws.Range('E43').Value = '  -0.35%  '
Set-TextCell 'D44' '8.179'
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 'D46' '1.005'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D47' '10.17'
$ws.Range('E47').Value = '  +0.44%  '
Set-TextCell 'D48' '102.47'
$ws.Range('E48').Value = '  -0.75%  '
Set-TextCell 'D49' '1.630'
$ws.Range('E49').Value = '  -0.36%  '
Set-TextCell 'D50' '0.06208'
$ws.Range('E50').Value = '  -0.18%  '
Set-TextCell 'D51' '64.86'
$ws.Range('E51').Value = '  +0.68%  '
